$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 21:55"

# Update Asturias row (row 20) figures
$ws.Range("B20").Value = 1384
$ws.Range("C20").Value = 135
$ws.Range("D20").Value = 1180
$ws.Range("E20").Value = 69
